$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7, 15, 6, 5),
    @(3, 8, 2, 12),
    @(5, 13, 3, 7),
    @(4, 8, 7, 12),
    @(4, 12, 5, 8),
    @(9, 14, 7, 6),
    @(4, 16, 5, 4),
    @(4, 13, 3, 7),
    @(5, 12, 6, 8),
    @(4, 5, 2, 15),
    @(4, 12, 3, 8),
    @(2, 12, 1, 8),
    @(3, 14, 5, 6),
    @(9, 3, 8, 17),
    @(5, 14, 6, 6),
    @(6, 19, 5, 1),
    @(7, 8, 6, 12),
    @(2, 14, 3, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 1454 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("A1472").Select()
$excel.ActiveWindow.ScrollRow = 1463
